$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 5 is:
#   "Spider challenge spider image: <url> " + "(free to use)"
# where the trailing space and "(free to use)" are two separate runs.
# Merge them into a single run's text " (free to use)" and drop the
# now-empty extra run.
$spaceRun = $tr.Characters(242, 1)
$spaceRun.Text = " (free to use)"

$oldParenRun = $tr.Characters(256, 13)
$oldParenRun.Text = ""

# Append a brand new paragraph right after paragraph 5 describing the
# landing-screen video credit/link.
$para5 = $tr.Paragraphs(5, 1)
$para5.InsertAfter("`rLanding screen video: https://www.pexels.com/video/lighted-candle-855262/") | Out-Null

$newPara = $tr.Paragraphs(6, 1)

# Split the new paragraph's text into two runs: "Landing screen video"
# and ": https://www.pexels.com/video/lighted-candle-855262/"
$urlStart = $newPara.Start + 20
$urlLen = $newPara.Length - 20
$urlRange = $tr.Characters($urlStart, $urlLen)
$urlRange.Text = ": https://www.pexels.com/video/lighted-candle-855262/"
